$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 72 was missing the formatting applied to the rest of the
#        "Descuentos no asociados a FC" / "MENORES VALORES" block (rows
#        66-71). Copy the formats (incl. the Text format on column G) down
#        from row 71 onto row 72 before touching its values, so the new
#        values pick up the right cell formatting/type.
$ws.Range("C71:I71").Copy() | Out-Null
$ws.Range("C72:I72").PasteSpecial(-4122) | Out-Null

# --- 2) "Motivo del descuento" (column G) codes were missing for the
#        RECHAZO rows 61-65; fill them in.
$ws.Range("G61").Value = "551"
$ws.Range("G62").Value = "551"
$ws.Range("G63").Value = "551"
$ws.Range("G64").Value = "551"
$ws.Range("G65").Value = "551"

# --- 3) Swap the WOB/384 discount-reason codes and flip the sign of the
#        matching invoice-amount / net-payment values for rows 66-72.
for ($r = 66; $r -le 72; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $gCell = $ws.Cells.Item($r, 7)
    $hCell = $ws.Cells.Item($r, 8)

    $eVal = $eCell.Value()
    $hVal = $hCell.Value()
    $gVal = $gCell.Value()

    $eCell.Value = -1 * $eVal
    $hCell.Value = -1 * $hVal

    if ($gVal -eq "WOB") {
        $gCell.Value = "384"
    } else {
        $gCell.Value = "WOB"
    }
}

# --- 4) Number-format adjustment: "Importe de factura" / "Pago Neto"
#        columns (E, H) move from 0.00 to #,##0.00 (thousands separator).
$ws.Range("E18").NumberFormat = "#,##0.00"
$ws.Range("H18").NumberFormat = "#,##0.00"
$ws.Range("E19:E72").NumberFormat = "#,##0.00"
$ws.Range("H19:H72").NumberFormat = "#,##0.00"
